$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 1.120168333333333
$ws.Range("H2").Value = 3.360505
$ws.Range("I2").Value = 0.001768092629909379
$ws.Range("J2").Value = 0.001768092629909379
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.229822
$ws.Range("N2").Value = 0.689466
$ws.Range("O2").Value = 0.09226175421862418
$ws.Range("P2").Value = 0.09226175421862419
$ws.Range("Q2").Value = 0.2574393267033333
$ws.Range("R2").Value = 2.31695394033
$ws.Range("S2").Value = 0.0001631273276564599
$ws.Range("T2").Value = 0.0001631273276564599
$ws.Range("G3").Value = 1.120168333333333
$ws.Range("H3").Value = 3.360505
$ws.Range("I3").Value = 0.001768092629909379
$ws.Range("J3").Value = 0.001768092629909379
$ws.Range("O3").Value = 0.4364142651333466
$ws.Range("P3").Value = 0.4364142651333466
$ws.Range("Q3").Value = 1.217733128219444
$ws.Range("R3").Value = 10.959598153975
$ws.Range("S3").Value = 0.0007716208457695876
$ws.Range("T3").Value = 0.0007716208457695876
$ws.Range("G4").Value = 1.120168333333333
$ws.Range("H4").Value = 3.360505
$ws.Range("I4").Value = 0.001768092629909379
$ws.Range("J4").Value = 0.001768092629909379
$ws.Range("M4").Value = 1.174057666666666
$ws.Range("O4").Value = 0.4713239806480292
$ws.Range("P4").Value = 0.4713239806480293
$ws.Range("Q4").Value = 1.315142219707222
$ws.Range("R4").Value = 11.836279977365
$ws.Range("S4").Value = 0.000833344456483331
$ws.Range("T4").Value = 0.0008333444564833312
$ws.Range("I5").Value = 0.9534130698726969
$ws.Range("J5").Value = 0.9534130698726969
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.229822
$ws.Range("N5").Value = 0.689466
$ws.Range("O5").Value = 0.09226175421862418
$ws.Range("P5").Value = 0.09226175421862419
$ws.Range("Q5").Value = 138.8196605914053
$ws.Range("R5").Value = 1249.376945322648
$ws.Range("S5").Value = 0.08796356232141872
$ws.Range("T5").Value = 0.08796356232141873
$ws.Range("I6").Value = 0.9534130698726969
$ws.Range("J6").Value = 0.9534130698726969
$ws.Range("O6").Value = 0.4364142651333466
$ws.Range("P6").Value = 0.4364142651333466
$ws.Range("Q6").Value = 656.6413209475844
$ws.Range("R6").Value = 5909.771888528259
$ws.Range("S6").Value = 0.416083064257021
$ws.Range("T6").Value = 0.416083064257021
$ws.Range("I7").Value = 0.9534130698726969
$ws.Range("J7").Value = 0.9534130698726969
$ws.Range("M7").Value = 1.174057666666666
$ws.Range("O7").Value = 0.4713239806480292
$ws.Range("P7").Value = 0.4713239806480293
$ws.Range("Q7").Value = 709.1674722237381
$ws.Range("S7").Value = 0.4493664432942571
$ws.Range("T7").Value = 0.4493664432942572
$ws.Range("H8").Value = 85.18441
$ws.Range("I8").Value = 0.04481883749739363
$ws.Range("J8").Value = 0.04481883749739363
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.229822
$ws.Range("N8").Value = 0.689466
$ws.Range("O8").Value = 0.09226175421862418
$ws.Range("P8").Value = 0.09226175421862419
$ws.Range("Q8").Value = 6.525750491673333
$ws.Range("R8").Value = 58.73175442506
$ws.Range("S8").Value = 0.004135064569548988
$ws.Range("T8").Value = 0.004135064569548989
$ws.Range("H9").Value = 85.18441
$ws.Range("I9").Value = 0.04481883749739363
$ws.Range("J9").Value = 0.04481883749739363
$ws.Range("O9").Value = 0.4364142651333466
$ws.Range("P9").Value = 0.4364142651333466
$ws.Range("R9").Value = 277.8114904109499
$ws.Range("S9").Value = 0.01955958003055592
$ws.Range("T9").Value = 0.01955958003055592
$ws.Range("H10").Value = 85.18441
$ws.Range("I10").Value = 0.04481883749739363
$ws.Range("J10").Value = 0.04481883749739363
$ws.Range("M10").Value = 1.174057666666666
$ws.Range("O10").Value = 0.4713239806480292
$ws.Range("P10").Value = 0.4713239806480293
$ws.Range("Q10").Value = 33.33713654699221
$ws.Range("S10").Value = 0.02112419289728872
$ws.Range("T10").Value = 0.02112419289728873
